$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet so the old 5x5 layout doesn't leave stray cells behind.
$ws.UsedRange.Clear()

# --- Row 1: report date header ---
$ws.Cells.Item(1,1).Value = "Date: 2025-06-19"

# --- Row 2: column headers ---
$ws.Cells.Item(2,1).Value = "Employee"
$ws.Cells.Item(2,2).Value = "Name"
$ws.Cells.Item(2,3).Value = "Location"
$ws.Cells.Item(2,4).Value = "Task"
$ws.Cells.Item(2,5).Value = "Clock In"
$ws.Cells.Item(2,6).Value = "Clock Out"
$ws.Cells.Item(2,7).Value = "Hours Worked"

# --- Rows 3-5: time entries ---
$ws.Cells.Item(3,1).Value = "olafur"
$ws.Cells.Item(3,2).Value = "Oli saer"
$ws.Cells.Item(3,3).Value = "Akureyri"
$ws.Cells.Item(3,4).Value = "Foundation prep"
$ws.Cells.Item(3,5).Value = "09:16"
$ws.Cells.Item(3,6).Value = "09:27"
$ws.Cells.Item(3,7).Value = 0.18

$ws.Cells.Item(4,1).Value = "olafur"
$ws.Cells.Item(4,2).Value = "Oli saer"
$ws.Cells.Item(4,3).Value = "Reykjavik"
$ws.Cells.Item(4,4).Value = "Pour concrete"
$ws.Cells.Item(4,5).Value = "09:45"
$ws.Cells.Item(4,6).Value = "09:49"
$ws.Cells.Item(4,7).Value = 0.07000000000000001

$ws.Cells.Item(5,1).Value = "olafur"
$ws.Cells.Item(5,2).Value = "Oli saer"
$ws.Cells.Item(5,3).Value = "Akureyri"
$ws.Cells.Item(5,4).Value = "Foundation prep"
$ws.Cells.Item(5,5).Value = "12:03"
$ws.Cells.Item(5,6).Value = "12:03"
$ws.Cells.Item(5,7).Value = 0

# --- Row 6: per-day total ---
$ws.Range("A6:F6").Font.Bold = $false
$ws.Cells.Item(6,7).Value = "Total: 0.25"

# --- Row 8: overall total ---
$ws.Range("A8:F8").Font.Bold = $false
$ws.Cells.Item(8,7).Value = "Overall Total Hours: 0.25"

# --- Row 9: spacer numeric cell (bold style, empty) ---
$ws.Cells.Item(9,7).Value = 0
$ws.Cells.Item(9,7).Font.Bold = $true
$ws.Cells.Item(9,7).ClearContents()

# --- Row 12: spacer numeric cell (bold style, empty) ---
$ws.Cells.Item(12,1).Value = 0
$ws.Cells.Item(12,1).Font.Bold = $true
$ws.Cells.Item(12,1).ClearContents()

# --- Row 13: task-summary header ---
$ws.Cells.Item(13,1).Value = "Task Name"
$ws.Cells.Item(13,2).Value = "Total Hours"

# --- Rows 14-15: task summary ---
$ws.Cells.Item(14,1).Value = "Foundation prep"
$ws.Cells.Item(14,2).Value = 0.18

$ws.Cells.Item(15,1).Value = "Pour concrete"
$ws.Cells.Item(15,2).Value = 0.07000000000000001

# --- Bold formatting for headline rows ---
$ws.Range("A1").Font.Bold = $true
$ws.Range("A2:G2").Font.Bold = $true
$ws.Range("A2:G2").HorizontalAlignment = -4108
$ws.Range("G6").Font.Bold = $true
$ws.Range("A13:B13").Font.Bold = $true

# --- New columns F and G get the same width treatment as the existing ones ---
$ws.Columns.Item(6).ColumnWidth = 13
$ws.Columns.Item(7).ColumnWidth = 13
